$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing column B (shifts old B.. onward right by 2)
$ws.Range("B1:C1").EntireColumn.Insert()

# New header cells (copy of A1's centered/wrap style is applied automatically by Insert)
$ws.Range("B1").Value = "temp resolution GURT"
$ws.Range("C1").Value = "temp resolution UTR2"

# Match column A's width on the two new columns
$ws.Range("B1:C1").EntireColumn.ColumnWidth = $ws.Range("A1").ColumnWidth

# Fill in the per-row time-resolution values for the data rows (2-14)
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 2).Value = 9.83
    $ws.Cells.Item($r, 3).Value = 7.94
}

# Highlight column A for the data rows with the new accent fill
$ws.Range("A2:A14").Interior.ThemeColor = 10

# Restore the view: scroll back to the frozen column and select the newly added column B data
$ws.Range("B2:B14").Select()
